# Auto-generated edit script: applies updated Coeurl Profits (cost/profit) figures
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 448.33334
$ws.Cells.Item(2, 9).Value = 600
$ws.Cells.Item(2, 11).Value = 600
$ws.Cells.Item(2, 13).Value = -487
$ws.Cells.Item(9, 8).Value = 129.125
$ws.Cells.Item(9, 9).Value = 146
$ws.Cells.Item(9, 11).Value = 146
$ws.Cells.Item(9, 13).Value = 23
$ws.Cells.Item(32, 8).Value = 2475
$ws.Cells.Item(32, 9).Value = 2500
$ws.Cells.Item(32, 10).Value = 2450
$ws.Cells.Item(32, 11).Value = 2500
$ws.Cells.Item(32, 12).Value = 2450
$ws.Cells.Item(32, 13).Value = -2174
$ws.Cells.Item(32, 14).Value = -3102
$ws.Cells.Item(98, 8).Value = 5108.727
$ws.Cells.Item(98, 10).Value = 8630
$ws.Cells.Item(98, 12).Value = 8630
$ws.Cells.Item(98, 14).Value = -11626
$ws.Cells.Item(112, 8).Value = 72721.8
$ws.Cells.Item(112, 9).Value = 3491
$ws.Cells.Item(112, 11).Value = 10473
$ws.Cells.Item(112, 13).Value = -9365
$ws.Cells.Item(122, 8).Value = 5108.727
$ws.Cells.Item(122, 10).Value = 8630
$ws.Cells.Item(122, 12).Value = 25890
$ws.Cells.Item(122, 14).Value = -30790
$ws.Cells.Item(134, 8).Value = 99999.75
$ws.Cells.Item(134, 10).Value = 99999.75
$ws.Cells.Item(134, 12).Value = 99999.75
$ws.Cells.Item(134, 14).Value = -110139.75

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10582.173
$ws.Cells.Item(32, 9).Value = 6358.595
$ws.Cells.Item(32, 10).Value = 28321.2
$ws.Cells.Item(32, 11).Value = 6358.595
$ws.Cells.Item(32, 12).Value = 28321.2
$ws.Cells.Item(32, 13).Value = -6071.595
$ws.Cells.Item(32, 14).Value = -28895.2
$ws.Cells.Item(37, 8).Value = 16166.667
$ws.Cells.Item(55, 8).Value = 19999
$ws.Cells.Item(55, 10).Value = 19999
$ws.Cells.Item(55, 12).Value = 19999
$ws.Cells.Item(55, 14).Value = -20629
$ws.Cells.Item(61, 8).Value = 3904.1765
$ws.Cells.Item(61, 9).Value = 2633.5833
$ws.Cells.Item(61, 11).Value = 2633.5833
$ws.Cells.Item(61, 13).Value = -2421.5833
$ws.Cells.Item(74, 8).Value = 5510.34
$ws.Cells.Item(74, 9).Value = 1746.4054
$ws.Cells.Item(74, 11).Value = 1746.4054
$ws.Cells.Item(74, 13).Value = -872.4054000000001
$ws.Cells.Item(77, 8).Value = 5510.34
$ws.Cells.Item(77, 9).Value = 1746.4054
$ws.Cells.Item(77, 11).Value = 8732.027
$ws.Cells.Item(77, 13).Value = -4364.027
$ws.Cells.Item(110, 8).Value = 4526.1353
$ws.Cells.Item(110, 9).Value = 6344.1113
$ws.Cells.Item(110, 10).Value = 2803.842
$ws.Cells.Item(110, 11).Value = 6344.1113
$ws.Cells.Item(110, 12).Value = 2803.842
$ws.Cells.Item(110, 13).Value = -4299.1113
$ws.Cells.Item(110, 14).Value = -6893.842000000001
$ws.Cells.Item(122, 8).Value = 12622.2
$ws.Cells.Item(122, 9).Value = 13277.75
$ws.Cells.Item(122, 10).Value = 10000
$ws.Cells.Item(122, 11).Value = 39833.25
$ws.Cells.Item(122, 12).Value = 30000
$ws.Cells.Item(122, 13).Value = -37383.25
$ws.Cells.Item(122, 14).Value = -34900
$ws.Cells.Item(136, 8).Value = 3904.1765
$ws.Cells.Item(136, 9).Value = 2633.5833
$ws.Cells.Item(136, 11).Value = 7900.749899999999
$ws.Cells.Item(136, 13).Value = -5350.749899999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 5923.5
$ws.Cells.Item(105, 9).Value = 6062
$ws.Cells.Item(105, 11).Value = 6062
$ws.Cells.Item(105, 13).Value = -4315
$ws.Cells.Item(134, 8).Value = 1364.7778
$ws.Cells.Item(134, 9).Value = 1340.3462
$ws.Cells.Item(134, 11).Value = 4021.0386
$ws.Cells.Item(134, 13).Value = -1486.0386

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 10928.143
$ws.Cells.Item(16, 9).Value = 8833.333000000001
$ws.Cells.Item(16, 11).Value = 8833.333000000001
$ws.Cells.Item(16, 13).Value = -8546.333000000001
$ws.Cells.Item(22, 8).Value = 331.66666
$ws.Cells.Item(22, 9).Value = 247.5
$ws.Cells.Item(22, 11).Value = 247.5
$ws.Cells.Item(22, 13).Value = 102.5
$ws.Cells.Item(58, 8).Value = 1892.9062
$ws.Cells.Item(58, 9).Value = 1388.3704
$ws.Cells.Item(58, 11).Value = 1388.3704
$ws.Cells.Item(58, 13).Value = -1185.3704
$ws.Cells.Item(86, 8).Value = 3499
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 3499
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 3499
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).Value = -5745
$ws.Cells.Item(87, 8).Value = 15914.75
$ws.Cells.Item(87, 10).Value = 15914.75
$ws.Cells.Item(87, 12).Value = 15914.75
$ws.Cells.Item(87, 14).Value = -18286.75
$ws.Cells.Item(89, 8).Value = 3499
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 3499
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 17495
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).Value = -28727
$ws.Cells.Item(90, 8).Value = 15914.75
$ws.Cells.Item(90, 10).Value = 15914.75
$ws.Cells.Item(90, 12).Value = 47744.25
$ws.Cells.Item(90, 14).Value = -59600.25
$ws.Cells.Item(113, 8).Value = 10928.143
$ws.Cells.Item(113, 9).Value = 8833.333000000001
$ws.Cells.Item(113, 11).Value = 8833.333000000001
$ws.Cells.Item(113, 13).Value = -6663.333000000001
$ws.Cells.Item(122, 8).Value = 1628.4667
$ws.Cells.Item(122, 9).Value = 1501.1666
$ws.Cells.Item(122, 10).Value = 2137.6667
$ws.Cells.Item(122, 11).Value = 4503.4998
$ws.Cells.Item(122, 12).Value = 6413.000100000001
$ws.Cells.Item(122, 13).Value = -2053.4998
$ws.Cells.Item(122, 14).Value = -11313.0001
$ws.Cells.Item(136, 8).Value = 1892.9062
$ws.Cells.Item(136, 9).Value = 1388.3704
$ws.Cells.Item(136, 11).Value = 4165.1112
$ws.Cells.Item(136, 13).Value = -1615.1112
$ws.Cells.Item(141, 8).Value = 139193.89
$ws.Cells.Item(141, 10).Value = 139193.89
$ws.Cells.Item(141, 12).Value = 139193.89
$ws.Cells.Item(141, 14).Value = -149553.89

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(48, 8).Value = 50
$ws.Cells.Item(48, 9).Value = 50
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 150
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 100
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(110, 8).Value = 22541.166
$ws.Cells.Item(110, 9).Value = 24599.4
$ws.Cells.Item(110, 10).Value = 12250
$ws.Cells.Item(110, 11).Value = 73798.20000000001
$ws.Cells.Item(110, 12).Value = 36750
$ws.Cells.Item(110, 13).Value = -69708.20000000001
$ws.Cells.Item(110, 14).Value = -44930
$ws.Cells.Item(116, 8).Value = 935
$ws.Cells.Item(116, 9).Value = 935
$ws.Cells.Item(116, 11).Value = 2805
$ws.Cells.Item(116, 13).Value = 637
$ws.Cells.Item(119, 8).Value = 999.5
$ws.Cells.Item(119, 9).Value = 999.5
$ws.Cells.Item(119, 11).Value = 2998.5
$ws.Cells.Item(119, 13).Value = 1839.5
$ws.Cells.Item(129, 8).Value = 1820.5555
$ws.Cells.Item(129, 9).Value = 794
$ws.Cells.Item(129, 11).Value = 2382
$ws.Cells.Item(129, 13).Value = 2618
$ws.Cells.Item(131, 8).Value = 3349.8262
$ws.Cells.Item(131, 9).Value = 1239.4445
$ws.Cells.Item(131, 10).Value = 3863.162
$ws.Cells.Item(131, 11).Value = 3718.3335
$ws.Cells.Item(131, 12).Value = 11589.486
$ws.Cells.Item(131, 13).Value = 1321.6665
$ws.Cells.Item(131, 14).Value = -21669.486

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 8).Value = 58000
$ws.Cells.Item(51, 10).Value = 58000
$ws.Cells.Item(51, 12).Value = 58000
$ws.Cells.Item(51, 14).Value = -59018
$ws.Cells.Item(57, 8).Value = 8344.727999999999
$ws.Cells.Item(70, 8).Value = 16164.228
$ws.Cells.Item(70, 9).Value = 15973.167
$ws.Cells.Item(70, 11).Value = 15973.167
$ws.Cells.Item(70, 13).Value = -15703.167
$ws.Cells.Item(73, 8).Value = 16164.228
$ws.Cells.Item(73, 9).Value = 15973.167
$ws.Cells.Item(73, 11).Value = 15973.167
$ws.Cells.Item(73, 13).Value = -15037.167
$ws.Cells.Item(80, 8).Value = 5124.125
$ws.Cells.Item(80, 9).Value = 3666
$ws.Cells.Item(80, 10).Value = 5999
$ws.Cells.Item(80, 11).Value = 3666
$ws.Cells.Item(80, 12).Value = 5999
$ws.Cells.Item(80, 13).Value = -2668
$ws.Cells.Item(80, 14).Value = -7995
$ws.Cells.Item(83, 8).Value = 5124.125
$ws.Cells.Item(83, 9).Value = 3666
$ws.Cells.Item(83, 10).Value = 5999
$ws.Cells.Item(83, 11).Value = 18330
$ws.Cells.Item(83, 12).Value = 29995
$ws.Cells.Item(83, 13).Value = -13338
$ws.Cells.Item(83, 14).Value = -39979
$ws.Cells.Item(102, 8).Value = 41668296
$ws.Cells.Item(102, 9).Value = 1506.5
$ws.Cells.Item(102, 10).Value = 250002260
$ws.Cells.Item(102, 11).Value = 1506.5
$ws.Cells.Item(102, 12).Value = 250002260
$ws.Cells.Item(102, 13).Value = 115.5
$ws.Cells.Item(102, 14).Value = -250005504
$ws.Cells.Item(122, 8).Value = 2088.2727
$ws.Cells.Item(122, 9).Value = 1997.1
$ws.Cells.Item(122, 11).Value = 5991.299999999999
$ws.Cells.Item(122, 13).Value = -3541.299999999999
$ws.Cells.Item(132, 8).Value = 7092.154
$ws.Cells.Item(132, 9).Value = 6024.75
$ws.Cells.Item(132, 11).Value = 18074.25
$ws.Cells.Item(132, 13).Value = -15544.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 40000
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 2219.8572
$ws.Cells.Item(46, 9).Value = 2066.6667
$ws.Cells.Item(46, 10).Value = 2334.75
$ws.Cells.Item(46, 11).Value = 2066.6667
$ws.Cells.Item(46, 12).Value = 2334.75
$ws.Cells.Item(46, 13).Value = -1878.6667
$ws.Cells.Item(46, 14).Value = -2710.75
$ws.Cells.Item(132, 8).Value = 4014.1667
$ws.Cells.Item(132, 9).Value = 3154.8333
$ws.Cells.Item(132, 11).Value = 9464.499899999999
$ws.Cells.Item(132, 13).Value = -6934.499899999999
$ws.Cells.Item(136, 8).Value = 5193.263
$ws.Cells.Item(136, 9).Value = 4829.7856
$ws.Cells.Item(136, 10).Value = 6211
$ws.Cells.Item(136, 11).Value = 14489.3568
$ws.Cells.Item(136, 12).Value = 18633
$ws.Cells.Item(136, 13).Value = -11939.3568
$ws.Cells.Item(136, 14).Value = -23733

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1473.3572
$ws.Cells.Item(113, 9).Value = 1379.2106
$ws.Cells.Item(113, 11).Value = 4137.6318
$ws.Cells.Item(113, 13).Value = -1967.6318
$ws.Cells.Item(122, 8).Value = 2400.276
$ws.Cells.Item(122, 10).Value = 2663.125
$ws.Cells.Item(122, 12).Value = 7989.375
$ws.Cells.Item(122, 14).Value = -12889.375
$ws.Cells.Item(135, 8).Value = 82437.5
$ws.Cells.Item(135, 9).Value = 35000
$ws.Cells.Item(135, 10).Value = 89214.28999999999
$ws.Cells.Item(135, 11).Value = 35000
$ws.Cells.Item(135, 12).Value = 89214.28999999999
$ws.Cells.Item(135, 13).Value = -29930
$ws.Cells.Item(135, 14).Value = -99354.28999999999
